$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '63.002.85'
$ws.Range("E2").Value = '  -2.05%  '

$ws.Range("D3").Value = '3.122.67'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  +0.03%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.53'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -2.80%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.56'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -4.36%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '3.116.74'
$ws.Range("E8").Value = '  -0.52%  '

$ws.Range("E9").Value = '  -1.95%  '

$ws.Range("E10").Value = '  -2.72%  '

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.29'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -0.69%  '

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  -2.95%  '

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  -2.38%  '

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.25'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("D15").Value = '3.634.74'
$ws.Range("E15").Value = '  -0.35%  '

$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").Value = '63.046.60'
$ws.Range("E17").Value = '  -1.84%  '

$ws.Range("D18").Value = '3.123.47'
$ws.Range("E18").Value = '  -0.31%  '

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.73'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  -1.58%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '477.32'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +0.31%  '

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.15'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -3.17%  '

$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.700'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -2.69%  '

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.69'
$ws.Range("D23").Style = $style

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.56'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  +3.09%  '

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.07'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -3.38%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("E27").Value = '  -2.61%  '

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.19'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -2.07%  '

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.00'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -5.69%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.09'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  +0.97%  '

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.11'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +1.79%  '

$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("E33").Value = '  -6.97%  '

$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.54'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  -3.01%  '

$ws.Range("E35").Value = '  -2.14%  '

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.84'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  -1.14%  '

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.99'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -0.94%  '

$ws.Range("D38").Value = '0.0₃0710'

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0390'
$ws.Range("D39").Style = $style

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '420.05'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -6.82%  '

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.116'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -0.97%  '

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.28'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  -0.31%  '

$ws.Range("E43").Value = '  -9.32%  '

$ws.Range("D44").Value = '2.882.49'
$ws.Range("E44").Value = '  +1.13%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.263'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -0.71%  '

$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.14'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -4.82%  '

$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  -0.06%  '

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.69'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -2.12%  '

$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("E50").Value = '  -5.40%  '

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.09'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -0.81%  '
